$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.912.57'
$ws.Range('E2').Value = '  +5.11%  '
$ws.Range('D3').Value = '2.347.84'
$ws.Range('E3').Value = '  +4.56%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').Value = '305.94'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '97.60'
$ws.Range('E6').Value = '  +3.06%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  +4.03%  '
$ws.Range('D10').Value = '35.78'
$ws.Range('E10').Value = '  +2.80%  '
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').Value = '7.42'
$ws.Range('E12').Value = '  +3.11%  '
$ws.Range('D14').Value = '2.709.97'
$ws.Range('E14').Value = '  +4.70%  '
$ws.Range('D15').Value = '2.345.02'
$ws.Range('E15').Value = '  +4.64%  '
$ws.Range('D16').Value = '14.20'
$ws.Range('E16').Value = '  +4.81%  '
$ws.Range('D17').Value = '0.829'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '46.793.68'
$ws.Range('E18').Value = '  +5.40%  '
$ws.Range('D19').Value = '13.65'
$ws.Range('E19').Value = '  +16.13%  '
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('D23').Value = '244.93'
$ws.Range('E23').Value = '  +3.14%  '
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').Value = '1.98'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').Value = '41.85'
$ws.Range('E27').Value = '  +13.46%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').Value = '9.87'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').Value = '20.15'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '5.76'
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('D32').Value = '152.44'
$ws.Range('E32').Value = '  +3.38%  '
$ws.Range('D33').Value = '0.0816'
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('D34').Value = '2.62'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').Value = '3.18'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('D38').Value = '1.82'
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('D39').Value = '4.03'
$ws.Range('E39').Value = '  +7.17%  '
$ws.Range('D40').Value = '0.0315'
$ws.Range('E40').Value = '  +5.18%  '
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('D42').Value = '13.79'
$ws.Range('E42').Value = '  -9.20%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.98'
$ws.Range('E44').Value = '  +12.70%  '
$ws.Range('D45').Value = '1.843.64'
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('E46').Value = '  +5.22%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '80.89'
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '73.72'
$ws.Range('E48').Value = '  +7.00%  '
$ws.Range('D49').Value = '4.94'
$ws.Range('E49').Value = '  +2.77%  '
$ws.Range('D50').Value = '98.54'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').Value = '55.16'
$ws.Range('E51').Value = '  +2.16%  '
